# The presentation ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> currently the "Office Theme" color scheme,
#                             used by the Notes Master (theme1 rel).
#   ppt/theme/theme2.xml  -> currently the "Integral" color scheme,
#                             used by the Slide Master (theme2 rel).
#
# The authored change swaps the two color schemes between the two theme
# parts (the Notes Master ends up with the "Integral" colors and the
# Slide Master ends up with the "Office Theme" colors), while the font
# scheme / format scheme (identical in both themes) and the two
# theme<->master relationships stay untouched.
#
# dk1/lt1 (black/white) are identical in both schemes, so only the other
# ten theme colors (dk2, lt2, accent1-6, hlink, folHlink) actually need
# to move.

function Hex2RGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Current "Office Theme" scheme (today on the Notes Master / theme1.xml).
$officeColors = @("44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

# Current "Integral" scheme (today on the Slide Master / theme2.xml).
$integralColors = @("455F51", "E3DED1", "99CB38", "63A537", "E6D024", "CC9700", "4EB3CF", "378DA6", "6B9F25", "B26B02")

$p = $ppt.ActivePresentation

# Slide Master's theme (theme2.xml): Integral -> Office Theme colors.
$slideTheme = $p.SlideMaster.Theme
$slideScheme = $slideTheme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $slideScheme.Item($i + 3).RGB = Hex2RGB $officeColors[$i]
}

# Notes Master's theme (theme1.xml): Office Theme -> Integral colors.
$notesTheme = $p.NotesMaster.Theme
$notesScheme = $notesTheme.ThemeColorScheme
for ($i = 0; $i -lt $integralColors.Count; $i++) {
    $notesScheme.Item($i + 3).RGB = Hex2RGB $integralColors[$i]
}
